$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to write so Excel does not
# reinterpret numeric-looking strings (e.g. "56.20", "9.00", "1.00") as numbers
# and strip significant trailing zeros / formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.228.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.588.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.584.47"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.55%  "
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.37%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.88"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.172.85"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.587.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.125"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.187.26"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.21"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.47"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.95"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.64"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.71"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "633.52"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.19"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.73"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.401"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0781"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.183.26"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.50%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.01"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0417"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.51%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.131"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.08"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.81"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.60"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.61%  "
